$d = $word.ActiveDocument

# Locate the insertion point: right at the start of the final (empty) paragraph,
# immediately after "This was an iterative process..." paragraph.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertRng = $lastPara.Range
$insertRng.End = $insertRng.Start

# Raw OOXML for the new "VI. Ethical Considerations" / "VII. Data Source Overview" /
# "VII. Project Plan" sections (32 paragraphs), followed by a blank trailing paragraph
# that restores the document-final empty paragraph InsertXML consumes.
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_Toc198289236"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>VI. Ethical Considerations</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">As with any project involving data, ethical considerations are of the utmost importance, especially when working with human subjects, even in anonymized form. Predictive </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>modeling</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> for human resource management raises issues of fairness, transparency, and biases that can be introduced by the data or </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>modeling</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> methods employed. A number of measures were taken to ensure that the project was following ethical best practices:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:bookmarkStart w:id="1" w:name="_Toc198289237"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>1. Data Privacy and Anonymization:</w:t></w:r><w:bookmarkEnd w:id="1"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>The data set used did not contain personally identifiable information (PII), thus complying with data protection regulations such as GDPR. There was no attempt to re-identify individuals.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:bookmarkStart w:id="2" w:name="_Toc198289238"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>2. Bias Detection and Mitigation:</w:t></w:r><w:bookmarkEnd w:id="2"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> As features like gender and city may be bias sources, their influence on model performance was carefully monitored. Models were evaluated for disparate </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:lastRenderedPageBreak/><w:t>impact, and techniques like balanced class weighting and fairness-aware preprocessing were explored.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:bookmarkStart w:id="3" w:name="_Toc198289239"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>3. Ethical Decision-Making:</w:t></w:r><w:bookmarkEnd w:id="3"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> The ultimate application of the model's predictions will not be utilized in punitive or discriminatory actions. Instead, they are to be utilized to inform positive retention efforts, such as by offering mentorship, flexible work schedules, or tailored engagement initiatives.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>In following these ethical guidelines, the project ensures responsible machine learning application in human resources.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="4" w:name="_Toc198289240"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>VII. Data Source Overview</w:t></w:r><w:bookmarkEnd w:id="4"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">The foundation of this predictive </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>modeling</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> project is established on the basis of a well-structured dataset that provides a detailed description of employee profiles. The dataset, available in the public domain, contains a range of features that span demographic and professional aspects. The main intention behind utilizing this dataset was to utilize pertinent attributes to make predictions about employee attrition.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>The dataset includes the following major features:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">• Education: </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Categorical variable indicating the highest educational qualification obtained (e.g., Bachelor's, Master's, PhD).</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">• City: </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Indicates geographical regions of employment (e.g., Bangalore, Pune, New Delhi).</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">• Gender: </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Binary or categorical data denoting the employee's gender.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">• Age: </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Numerical variable denoting the age of the employee.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">• </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>ExperienceInCurrentDomain</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Numerical value indicating the number of years the employee has been in the current line of business.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">• </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>EverBenched</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>A dummy variable indicating if the employee was ever on the bench (i.e., briefly without a work assignment).</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">• </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>LeaveOrNot</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>The response variable indicating if the employee stayed (0) or left (1) the company.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Although not exhaustive with regard to all possible employment factors, the dataset provides a diverse and representative representation that allows for predictive analytics and data-driven HR interventions.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">Initial data was provided in CSV format, which was imported into a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Jupyter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> Notebook platform for exploratory data analysis and preprocessing. Minimal integrity checks for data preparation were carried out to ensure that the data was ready and consistent for model build. Missing values, data type issues, and category consistency issues were handled in the preprocessing step.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="5" w:name="_Toc198289241"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>VII. Project Plan</w:t></w:r><w:bookmarkEnd w:id="5"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:lastRenderedPageBreak/><w:t>The project plan adhered to a logical sequence as per CRISP-DM, with checkpoints and iterative evaluation in between. The following steps were formulated and executed:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:bookmarkStart w:id="6" w:name="_Toc198289242"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>1. Requirement Gathering and Problem Definition:</w:t></w:r><w:bookmarkEnd w:id="6"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> Established business goals and translated them into data mining objectives.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:bookmarkStart w:id="7" w:name="_Toc198289243"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>2. Data Understanding and Preprocessing:</w:t></w:r><w:bookmarkEnd w:id="7"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> Performed intensive exploratory analysis, identified missing or inconsistent data, and used appropriate transformation techniques.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:bookmarkStart w:id="8" w:name="_Toc198289244"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>3. Feature Engineering:</w:t></w:r><w:bookmarkEnd w:id="8"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> Derived new features and carried out feature selection based on correlation and significance.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:bookmarkStart w:id="9" w:name="_Toc198289245"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>4. Model Building:</w:t></w:r><w:bookmarkEnd w:id="9"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> Several machine learning models were tried, compared, and assessed.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:bookmarkStart w:id="10" w:name="_Toc198289246"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>5. Model Tuning:</w:t></w:r><w:bookmarkEnd w:id="10"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> Hyperparameters were optimized for improved model performance.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:bookmarkStart w:id="11" w:name="_Toc198289247"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>6. Interpretation of Results</w:t></w:r><w:bookmarkEnd w:id="11"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> Visualized results, outlined model results, and explained the impact of features.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:bookmarkStart w:id="12" w:name="_Toc198289248"/><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>7. Reporting and Documentation:</w:t></w:r><w:bookmarkEnd w:id="12"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> Presented findings in a report and summarized results graphically for poster presentation. This structured workflow helped ensure that all aspects of the project were thoroughly addressed, from data wrangling to strategic insights.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p><w:p/>
'@
$insertRng.InsertXML($xml)

# InsertXML (this host) drops <w:rStyle> from inserted runs, so the numbered
# sub-heading runs ("1. Data Privacy and Anonymization:", etc.) need the
# "Heading2Char" character style re-applied after the fact via Find + Range.Style.
$headings = @(
    '1. Data Privacy and Anonymization:',
    '2. Bias Detection and Mitigation:',
    '3. Ethical Decision-Making:',
    '1. Requirement Gathering and Problem Definition:',
    '2. Data Understanding and Preprocessing:',
    '3. Feature Engineering:',
    '4. Model Building:',
    '5. Model Tuning:',
    '6. Interpretation of Results',
    '7. Reporting and Documentation:',
)

foreach ($h in $headings) {
    $fr = $d.Content
    $ok = $fr.Find.Execute($h, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $fr.Style = "Heading2Char"
    } else {
        Write-Output ("NOT FOUND: " + $h)
    }
}

# The heading "1. Data Privacy and Anonymization:" is followed by a separate run
# containing just a space, which also carries the Heading2Char style in the source.
$fr2 = $d.Content
$ok2 = $fr2.Find.Execute("1. Data Privacy and Anonymization:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok2) {
    $fr2.End = $fr2.End + 1
    $fr2.Style = "Heading2Char"
}

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
